# Update the "timestamp" column (Z) on the active sheet to reflect the
# latest re-run of the pcsmote sampling log. Every data row (2-112) gets a
# refreshed timestamp; rows generated in the same batch share the exact
# same microsecond-precision value, matching the ranges below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z45").Value   = "2025-10-17T07:09:37.141529"
$ws.Range("Z46:Z52").Value  = "2025-10-17T07:09:37.239864"
$ws.Range("Z53:Z56").Value  = "2025-10-17T07:09:37.240863"
$ws.Range("Z57:Z64").Value  = "2025-10-17T07:09:37.241866"
$ws.Range("Z65:Z74").Value  = "2025-10-17T07:09:37.242865"
$ws.Range("Z75:Z78").Value  = "2025-10-17T07:09:37.332229"
$ws.Range("Z79").Value      = "2025-10-17T07:09:37.340814"
$ws.Range("Z80:Z102").Value = "2025-10-17T07:09:37.341543"
$ws.Range("Z103:Z104").Value = "2025-10-17T07:09:37.424822"
$ws.Range("Z105:Z112").Value = "2025-10-17T07:09:37.429836"
